$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 1: "100" -> "0M"
$tbl.Cell(1, 1).Range.Text = "0M"
# Row 2: "0" -> "0M"
$tbl.Cell(2, 1).Range.Text = "0M"
# Row 3: "416" -> "0M"
$tbl.Cell(3, 1).Range.Text = "0M"
# Row 4: "3" -> "44"
$tbl.Cell(4, 1).Range.Text = "44"
# Row 5: "0.00003" -> "0.00002"
$tbl.Cell(5, 1).Range.Text = "0.00002"
# Row 12: "0.00012" -> "0.00156"
$tbl.Cell(12, 1).Range.Text = "0.00156"

# Rows 44, 45, 46: collapse the multi-run tab-separated text down to a single value
$tbl.Cell(44, 1).Range.Text = "100"
$tbl.Cell(45, 1).Range.Text = "0"
$tbl.Cell(46, 1).Range.Text = "416"
